# Update "想去人数" (people interested) counts in the F column across
# sheets "展览", "演出" and "全部类型" to reflect freshly generated data
# (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1707
$ws.Range("F3").Value = 9195
$ws.Range("F7").Value = 1398
$ws.Range("F11").Value = 5996
$ws.Range("F15").Value = 4740
$ws.Range("F17").Value = 169
$ws.Range("F19").Value = 40
$ws.Range("F20").Value = 347
$ws.Range("F21").Value = 37
$ws.Range("F23").Value = 271
$ws.Range("F25").Value = 3182
$ws.Range("F26").Value = 133

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 54

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1707
$ws.Range("F3").Value = 9195
$ws.Range("F5").Value = 54
$ws.Range("F8").Value = 1398
$ws.Range("F12").Value = 5996
$ws.Range("F16").Value = 4740
$ws.Range("F18").Value = 169
$ws.Range("F20").Value = 40
$ws.Range("F21").Value = 347
$ws.Range("F22").Value = 37
$ws.Range("F24").Value = 271
$ws.Range("F26").Value = 3182
$ws.Range("F28").Value = 133

$wb.Save()
